# Updates Price (D) and Volume(1h) (E) columns for the cryptos list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.649.60"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "2.801.06"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0840"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.13%  "
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").Value = "3.241.87"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Value = "2.796.01"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "51.655.91"
$ws.Range("E19").Value = "  +3.43%  "
$ws.Range("E20").Value = "  +3.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.20%  "
$ws.Range("E31").Value = "  +4.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "52.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0448"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0856"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  -5.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("D46").Value = "2.140.03"
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.85%  "
$ws.Range("E48").Value = "  +7.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.924"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.18%  "
$ws.Range("E50").Value = "  +10.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.220"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +15.83%  "
